$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set cells whose new text would otherwise be parsed as numbers to Text format first,
# so Excel keeps the original textual value (matching the source inline strings).

$ws.Range("D2").Value = "25.894.27"
$ws.Range("E2").Value = "  -0.18%  "
$ws.Range("D3").Value = "1.584.23"
$ws.Range("E3").Value = "  -1.95%  "
$ws.Range("E4").Value = "  -0.28%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.23"
$ws.Range("E5").Value = "  -0.54%  "
$ws.Range("E6").Value = "  -0.29%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.477"
$ws.Range("E7").Value = "  -2.02%  "
$ws.Range("E8").Value = "  -0.13%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0614"
$ws.Range("E9").Value = "  -0.83%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.12"
$ws.Range("E10").Value = "  -0.16%  "
$ws.Range("E11").Value = "  -0.14%  "
$ws.Range("D12").Value = "1.803.49"
$ws.Range("D13").Value = "1.579.83"
$ws.Range("E13").Value = "  -2.10%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.03"
$ws.Range("E14").Value = "  -1.30%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.504"
$ws.Range("E15").Value = "  -2.66%  "
$ws.Range("D16").Value = "25.873.80"
$ws.Range("E16").Value = "  -0.30%  "
$ws.Range("D17").Value = "0.0₃0723"
$ws.Range("E17").Value = "  -1.00%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "60.06"
$ws.Range("E18").Value = "  -2.33%  "
$ws.Range("E19").Value = "  -0.17%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "193.29"
$ws.Range("E20").Value = "  +1.15%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.19"
$ws.Range("E21").Value = "  -0.73%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.35"
$ws.Range("E22").Value = "  -0.63%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.93"
$ws.Range("E23").Value = "  -1.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.131"
$ws.Range("E24").Value = "  +0.91%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "140.76"
$ws.Range("E25").Value = "  -1.51%  "
$ws.Range("E26").Value = "  -0.17%  "
$ws.Range("E27").Value = "  -1.54%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.07"
$ws.Range("E28").Value = "  -0.20%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.44"
$ws.Range("E29").Value = "  -2.46%  "
$ws.Range("E30").Value = "  -4.46%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0472"
$ws.Range("E31").Value = "  -0.03%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.13"
$ws.Range("E32").Value = "  +0.69%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.02"
$ws.Range("E33").Value = "  -1.90%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.50"
$ws.Range("E34").Value = "  +0.94%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.35"
$ws.Range("E35").Value = "  -2.07%  "
$ws.Range("D36").Value = "1.096.37"
$ws.Range("E36").Value = "  -2.24%  "
$ws.Range("E37").Value = "  -0.35%  "
$ws.Range("E38").Value = "  -1.82%  "
$ws.Range("E39").Value = "  -0.71%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.501"
$ws.Range("E40").Value = "  -1.75%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.778"
$ws.Range("E41").Value = "  -4.95%  "
$ws.Range("E42").Value = "  +6.95%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "93.12"
$ws.Range("E43").Value = "  -4.03%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.09"
$ws.Range("E44").Value = "  +0.50%  "
$ws.Range("D45").Value = "1.717.67"
$ws.Range("E45").Value = "  -1.93%  "
$ws.Range("D46").Value = "0.0₆0108"
$ws.Range("E46").Value = "  -5.34%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.51"
$ws.Range("E47").Value = "  +1.82%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "53.12"
$ws.Range("E48").Value = "  -1.10%  "
$ws.Range("E49").Value = "  -1.48%  "
$ws.Range("E50").Value = "  -0.94%  "
$ws.Range("E51").Value = "  -0.21%  "
